$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "technical"
$ws.Range("B7").Value = 1

$ws.Range("A8").Value = "back"
$ws.Range("B8").Value = 2
